$wb = $excel.ActiveWorkbook

$wsExtraction = $wb.Worksheets.Item("Extraction")
$wsExtraction.Rows(8).Delete()

$wsOverlap = $wb.Worksheets.Item("Overlap")
$wsOverlap.Rows(8).Delete()
